$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4211.1113
$ws.Range("I40").Value = 1975
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 1975
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -1800
$ws.Range("N40").Value = -6350

$ws.Range("H64").Value = 3145.389
$ws.Range("I64").Value = 3161.6
$ws.Range("J64").Value = 3139.1538
$ws.Range("K64").Value = 3161.6
$ws.Range("L64").Value = 3139.1538
$ws.Range("M64").Value = -2913.6
$ws.Range("N64").Value = -3635.1538

$ws.Range("H67").Value = 3145.389
$ws.Range("I67").Value = 3161.6
$ws.Range("J67").Value = 3139.1538
$ws.Range("K67").Value = 3161.6
$ws.Range("L67").Value = 3139.1538
$ws.Range("M67").Value = -2303.6
$ws.Range("N67").Value = -4855.1538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 817
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 1125.5
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 1125.5
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -1349.5

$ws.Range("H32").Value = 9863.691999999999
$ws.Range("I32").Value = 3308.3555
$ws.Range("K32").Value = 3308.3555
$ws.Range("M32").Value = -3021.3555

$ws.Range("H63").Value = 4011.111
$ws.Range("I63").Value = 3662.5
$ws.Range("K63").Value = 3662.5
$ws.Range("M63").Value = -2976.5

$ws.Range("H66").Value = 4011.111
$ws.Range("I66").Value = 3662.5
$ws.Range("K66").Value = 18312.5
$ws.Range("M66").Value = -14880.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 817
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 1125.5
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 1125.5
$ws.Range("M4").Value = -85
$ws.Range("N4").Value = -1355.5

$ws.Range("H82").Value = 9607.1
$ws.Range("I82").Value = 6841.222
$ws.Range("J82").Value = 34500
$ws.Range("K82").Value = 6841.222
$ws.Range("L82").Value = 34500
$ws.Range("M82").Value = -6458.222
$ws.Range("N82").Value = -35266

$ws.Range("H85").Value = 9607.1
$ws.Range("I85").Value = 6841.222
$ws.Range("J85").Value = 34500
$ws.Range("K85").Value = 6841.222
$ws.Range("L85").Value = 34500
$ws.Range("M85").Value = -5515.222
$ws.Range("N85").Value = -37152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 104.28571
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 190
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 190
$ws.Range("M7").Value = 73
$ws.Range("N7").Value = -416

$ws.Range("H31").Value = 7019216.5
$ws.Range("I31").Value = 8334489
$ws.Range("J31").Value = 6062654.5
$ws.Range("K31").Value = 8334489
$ws.Range("L31").Value = 6062654.5
$ws.Range("M31").Value = -8334194
$ws.Range("N31").Value = -6063244.5

$ws.Range("H34").Value = 7019216.5
$ws.Range("I34").Value = 8334489
$ws.Range("J34").Value = 6062654.5
$ws.Range("K34").Value = 8334489
$ws.Range("L34").Value = 6062654.5
$ws.Range("M34").Value = -8334287
$ws.Range("N34").Value = -6063058.5

$ws.Range("H41").Value = 16353
$ws.Range("I41").Value = 950
$ws.Range("J41").Value = 26621.666
$ws.Range("K41").Value = 950
$ws.Range("L41").Value = 26621.666
$ws.Range("M41").Value = -522
$ws.Range("N41").Value = -27477.666

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = 0

$ws.Range("H51").Value = 24000
$ws.Range("J51").Value = 24000
$ws.Range("L51").Value = 24000
$ws.Range("N51").Value = -25472

$ws.Range("H59").Value = 19395
$ws.Range("J59").Value = 19395
$ws.Range("L59").Value = 19395
$ws.Range("N59").Value = -21685

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = 0

$ws.Range("H61").Value = 24000
$ws.Range("J61").Value = 24000
$ws.Range("L61").Value = 24000
$ws.Range("N61").Value = -24696

$ws.Range("H62").Value = 111114184
$ws.Range("I62").Value = 3460
$ws.Range("J62").Value = 250002600
$ws.Range("K62").Value = 3460
$ws.Range("L62").Value = 250002600
$ws.Range("M62").Value = -2836
$ws.Range("N62").Value = -250003848

$ws.Range("H65").Value = 111114184
$ws.Range("I65").Value = 3460
$ws.Range("J65").Value = 250002600
$ws.Range("K65").Value = 17300
$ws.Range("L65").Value = 1250013000
$ws.Range("M65").Value = -14180
$ws.Range("N65").Value = -1250019240

$ws.Range("H68").Value = 10268
$ws.Range("I68").Value = 10268
$ws.Range("K68").Value = 10268
$ws.Range("M68").Value = -9519

$ws.Range("H71").Value = 10268
$ws.Range("I71").Value = 10268
$ws.Range("K71").Value = 30804
$ws.Range("M71").Value = -27060

$ws.Range("H74").Value = 29314
$ws.Range("J74").Value = 29314
$ws.Range("L74").Value = 29314
$ws.Range("N74").Value = -31062

$ws.Range("H77").Value = 29314
$ws.Range("J77").Value = 29314
$ws.Range("L77").Value = 87942
$ws.Range("N77").Value = -96678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1883.3334
$ws.Range("I46").Value = 1600
$ws.Range("J46").Value = 2166.6667
$ws.Range("K46").Value = 1600
$ws.Range("L46").Value = 2166.6667
$ws.Range("M46").Value = -1412
$ws.Range("N46").Value = -2542.6667

$ws.Range("H68").Value = 11668696
$ws.Range("I68").Value = 22557222
$ws.Range("J68").Value = 2418.7144
$ws.Range("K68").Value = 22557222
$ws.Range("L68").Value = 2418.7144
$ws.Range("M68").Value = -22556473
$ws.Range("N68").Value = -3916.7144

$ws.Range("H71").Value = 11668696
$ws.Range("I71").Value = 22557222
$ws.Range("J71").Value = 2418.7144
$ws.Range("K71").Value = 112786110
$ws.Range("L71").Value = 12093.572
$ws.Range("M71").Value = -112782366
$ws.Range("N71").Value = -19581.572
